$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row update: insert Jan_2026 before Dec_2025, drop Oct_2025 (old col F)
$ws.Cells.Item(1,4).Value = "Jan_2026"
$ws.Cells.Item(1,5).Value = "Dec_2025"
$ws.Cells.Item(1,6).Value = "Nov_2025"

# Data rows (2..31) reflect the resorted/refreshed holdings table
$ws.Cells.Item(2,1).Value = "INE406A01037"
$ws.Cells.Item(2,2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(2,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(2,4).Value = 7.921843
$ws.Cells.Item(2,5).Value = 6.98705
$ws.Cells.Item(2,6).Value = 6.962837
$ws.Cells.Item(2,7).Value = 0.934793
$ws.Cells.Item(2,8).Value = 0.9590059999999996

$ws.Cells.Item(3,1).Value = "INE775A01035"
$ws.Cells.Item(3,2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(3,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(3,4).Value = 6.56115
$ws.Cells.Item(3,5).Value = 6.276131
$ws.Cells.Item(3,6).Value = 5.439641
$ws.Cells.Item(3,7).Value = 0.2850189999999992
$ws.Cells.Item(3,8).Value = 1.121509

$ws.Cells.Item(4,1).Value = "INE281B01032"
$ws.Cells.Item(4,2).Value = "Lloyds Metals And Energy Limited"
$ws.Cells.Item(4,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(4,4).Value = 6.354931
$ws.Cells.Item(4,5).Value = 6.849578
$ws.Cells.Item(4,6).Value = 6.079997
$ws.Cells.Item(4,7).Value = -0.4946470000000005
$ws.Cells.Item(4,8).Value = 0.274934

$ws.Cells.Item(5,1).Value = "INE018A01030"
$ws.Cells.Item(5,2).Value = "Larsen & Toubro Limited"
$ws.Cells.Item(5,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(5,4).Value = 5.980584
$ws.Cells.Item(5,5).Value = 5.592055
$ws.Cells.Item(5,6).Value = 5.355863
$ws.Cells.Item(5,7).Value = 0.3885290000000001
$ws.Cells.Item(5,8).Value = 0.6247210000000001

$ws.Cells.Item(6,1).Value = "INE002A01018"
$ws.Cells.Item(6,2).Value = "Reliance Industries Limited"
$ws.Cells.Item(6,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(6,4).Value = 5.877455
$ws.Cells.Item(6,5).Value = 10.008666
$ws.Cells.Item(6,6).Value = 10.111124
$ws.Cells.Item(6,7).Value = -4.131211
$ws.Cells.Item(6,8).Value = -4.233669

$ws.Cells.Item(7,1).Value = "INE814H01029"
$ws.Cells.Item(7,2).Value = "Adani Power Limited"
$ws.Cells.Item(7,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(7,4).Value = 5.445909
$ws.Cells.Item(7,5).Value = 5.17921
$ws.Cells.Item(7,6).Value = 5.13404
$ws.Cells.Item(7,7).Value = 0.266699
$ws.Cells.Item(7,8).Value = 0.3118690000000006

$ws.Cells.Item(8,1).Value = "INE758E01017"
$ws.Cells.Item(8,2).Value = "Jio Financial Services Limited"
$ws.Cells.Item(8,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(8,4).Value = 3.638142
$ws.Cells.Item(8,5).Value = 3.796489
$ws.Cells.Item(8,6).Value = 3.787719
$ws.Cells.Item(8,7).Value = -0.1583469999999996
$ws.Cells.Item(8,8).Value = -0.1495769999999998

$ws.Cells.Item(9,1).Value = "INE795G01014"
$ws.Cells.Item(9,2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(9,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(9,4).Value = 2.919895
$ws.Cells.Item(9,5).Value = 2.696909
$ws.Cells.Item(9,6).Value = 2.641767
$ws.Cells.Item(9,7).Value = 0.2229860000000001
$ws.Cells.Item(9,8).Value = 0.2781279999999997

$ws.Cells.Item(10,1).Value = "INE930H01031"
$ws.Cells.Item(10,2).Value = "K.P.R. Mill Limited"
$ws.Cells.Item(10,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(10,4).Value = 2.616977
$ws.Cells.Item(10,5).Value = 2.553775
$ws.Cells.Item(10,6).Value = 2.81131
$ws.Cells.Item(10,7).Value = 0.06320199999999998
$ws.Cells.Item(10,8).Value = -0.1943330000000003

$ws.Cells.Item(11,1).Value = "INE216A01030"
$ws.Cells.Item(11,2).Value = "Britannia Industries Limited"
$ws.Cells.Item(11,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(11,4).Value = 2.478869
$ws.Cells.Item(11,5).Value = 3.216406
$ws.Cells.Item(11,6).Value = 2.996257
$ws.Cells.Item(11,7).Value = -0.7375370000000001
$ws.Cells.Item(11,8).Value = -0.517388

$ws.Cells.Item(12,1).Value = "INE151A01013"
$ws.Cells.Item(12,2).Value = "Tata Communications Limited"
$ws.Cells.Item(12,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(12,4).Value = 2.22615
$ws.Cells.Item(12,5).Value = 2.329661
$ws.Cells.Item(12,6).Value = 2.226495
$ws.Cells.Item(12,7).Value = -0.1035110000000001
$ws.Cells.Item(12,8).Value = -0.0003449999999998177

$ws.Cells.Item(13,1).Value = "INE699H01024"
$ws.Cells.Item(13,2).Value = "Adani Wilmar Limited"
$ws.Cells.Item(13,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(13,4).Value = 2.016296
$ws.Cells.Item(13,5).Value = 2.013451
$ws.Cells.Item(13,6).Value = 2.111941
$ws.Cells.Item(13,7).Value = 0.002845000000000208
$ws.Cells.Item(13,8).Value = -0.09564499999999976

$ws.Cells.Item(14,1).Value = "INE042A01014"
$ws.Cells.Item(14,2).Value = "Escorts Kubota Limited"
$ws.Cells.Item(14,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(14,4).Value = 1.933665
$ws.Cells.Item(14,5).Value = 2.099768
$ws.Cells.Item(14,6).Value = 2.070896
$ws.Cells.Item(14,7).Value = -0.1661030000000001
$ws.Cells.Item(14,8).Value = -0.1372309999999999

$ws.Cells.Item(15,1).Value = "INE881D01027"
$ws.Cells.Item(15,2).Value = "Oracle Financial Services Software Ltd"
$ws.Cells.Item(15,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(15,4).Value = 1.637482
$ws.Cells.Item(15,5).Value = 1.46054
$ws.Cells.Item(15,6).Value = 1.48032
$ws.Cells.Item(15,7).Value = 0.1769420000000002
$ws.Cells.Item(15,8).Value = 0.157162

$ws.Cells.Item(16,1).Value = "INE364U01010"
$ws.Cells.Item(16,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(16,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(16,4).Value = 1.60878
$ws.Cells.Item(16,5).Value = 1.725262
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = -0.116482
$ws.Cells.Item(16,8).Value = 1.60878

$ws.Cells.Item(17,1).Value = "INE470A01017"
$ws.Cells.Item(17,2).Value = "3M India Limited"
$ws.Cells.Item(17,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(17,4).Value = 1.380572
$ws.Cells.Item(17,5).Value = 1.269292
$ws.Cells.Item(17,6).Value = 1.21028
$ws.Cells.Item(17,7).Value = 0.1112799999999998
$ws.Cells.Item(17,8).Value = 0.1702919999999999

$ws.Cells.Item(18,1).Value = "INE880J01026"
$ws.Cells.Item(18,2).Value = "JSW Infrastructure Limited"
$ws.Cells.Item(18,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(18,4).Value = 1.295972
$ws.Cells.Item(18,5).Value = 1.285633
$ws.Cells.Item(18,6).Value = 1.172849
$ws.Cells.Item(18,7).Value = 0.01033899999999988
$ws.Cells.Item(18,8).Value = 0.1231229999999999

$ws.Cells.Item(19,1).Value = "INE0BS701011"
$ws.Cells.Item(19,2).Value = "Premier Energies Limited"
$ws.Cells.Item(19,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(19,4).Value = 1.184429
$ws.Cells.Item(19,5).Value = 1.24758
$ws.Cells.Item(19,6).Value = 1.389102
$ws.Cells.Item(19,7).Value = -0.06315099999999996
$ws.Cells.Item(19,8).Value = -0.2046730000000001

$ws.Cells.Item(20,1).Value = "INE0J1Y01017"
$ws.Cells.Item(20,2).Value = "Life Insurance Corporation Of India"
$ws.Cells.Item(20,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(20,4).Value = 1.150905
$ws.Cells.Item(20,5).Value = 4.069492
$ws.Cells.Item(20,6).Value = 4.09185
$ws.Cells.Item(20,7).Value = -2.918587
$ws.Cells.Item(20,8).Value = -2.940945

$ws.Cells.Item(21,1).Value = "INE931S01010"
$ws.Cells.Item(21,2).Value = "Adani Energy Solutions Limited"
$ws.Cells.Item(21,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(21,4).Value = 1.107085
$ws.Cells.Item(21,5).Value = 1.144499
$ws.Cells.Item(21,6).Value = 1.064786
$ws.Cells.Item(21,7).Value = -0.03741399999999984
$ws.Cells.Item(21,8).Value = 0.04229900000000009

$ws.Cells.Item(22,1).Value = "INE259A01022"
$ws.Cells.Item(22,2).Value = "Colgate-Palmolive (India) Ltd"
$ws.Cells.Item(22,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(22,4).Value = 0.989349
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 0.989349
$ws.Cells.Item(22,8).Value = 0.989349

$ws.Cells.Item(23,1).Value = "INE179A01014"
$ws.Cells.Item(23,2).Value = "Procter & Gamble Hygiene & Health Care Limited"
$ws.Cells.Item(23,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(23,4).Value = 0.534404
$ws.Cells.Item(23,5).Value = 0.530464
$ws.Cells.Item(23,6).Value = 0.501344
$ws.Cells.Item(23,7).Value = 0.003939999999999944
$ws.Cells.Item(23,8).Value = 0.03305999999999998

$ws.Cells.Item(24,1).Value = "INE018E01016"
$ws.Cells.Item(24,2).Value = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(24,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(24,4).Value = 0.415008
$ws.Cells.Item(24,5).Value = 0.42731
$ws.Cells.Item(24,6).Value = 0.419452
$ws.Cells.Item(24,7).Value = -0.01230200000000004
$ws.Cells.Item(24,8).Value = -0.004444000000000004

$ws.Cells.Item(25,1).Value = "INE271C01023"
$ws.Cells.Item(25,2).Value = "DLF Limited"
$ws.Cells.Item(25,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 3.552037
$ws.Cells.Item(25,6).Value = 3.593398
$ws.Cells.Item(25,7).Value = -3.552037
$ws.Cells.Item(25,8).Value = -3.593398

$ws.Cells.Item(26,1).Value = "INE467B01029"
$ws.Cells.Item(26,2).Value = "Tata Consultancy Services Limited"
$ws.Cells.Item(26,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(26,4).Value = 0
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 2.952062
$ws.Cells.Item(26,7).Value = 0
$ws.Cells.Item(26,8).Value = -2.952062

$ws.Cells.Item(27,1).Value = "INE686F01025"
$ws.Cells.Item(27,2).Value = "UNITED BREWERIES LIMITED"
$ws.Cells.Item(27,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(27,4).Value = 0
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 0
$ws.Cells.Item(27,7).Value = 0
$ws.Cells.Item(27,8).Value = -0.443915

$ws.Cells.Item(28,1).Value = "INE154A01025"
$ws.Cells.Item(28,2).Value = "ITC Limited"
$ws.Cells.Item(28,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 3.439933
$ws.Cells.Item(28,7).Value = -3.439933
$ws.Cells.Item(28,8).Value = -3.316147

$ws.Cells.Item(29,1).Value = "INE14LE01019"
$ws.Cells.Item(29,2).Value = "Aditya Birla Lifestyle Brands Limited"
$ws.Cells.Item(29,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = 0.683092
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = -0.683092

$ws.Cells.Item(30,1).Value = "INE776C01039"
$ws.Cells.Item(30,2).Value = "GMR Airports Limited"
$ws.Cells.Item(30,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 2.514308
$ws.Cells.Item(30,7).Value = -2.514308
$ws.Cells.Item(30,8).Value = -3.026149

$ws.Cells.Item(31,1).Value = "INE115A01026"
$ws.Cells.Item(31,2).Value = "LIC Housing Finance Ltd"
$ws.Cells.Item(31,3).Value = "quant Large and Mid Cap Fund"
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = 0
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(31,8).Value = -0.146427

# Remove the two now-unused trailing rows (table shrank from 32 to 30 data rows)
$ws.Rows.Item(33).Delete()
$ws.Rows.Item(32).Delete()
